# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Only column G (header "K") changes, rows 2-34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..34 (1:1 with existing row order / dates already on sheet)
$newK = @{
    2  = 6
    3  = 3
    4  = 6
    5  = 7
    6  = 4
    7  = 4
    8  = 6
    9  = 8
    10 = 7
    11 = 5
    12 = 4
    13 = 5
    14 = 8
    15 = 4
    16 = 5
    17 = 2
    18 = 4
    19 = 4
    20 = 8
    21 = 7
    22 = 7
    23 = 5
    24 = 5
    25 = 6
    26 = 9
    27 = 4
    28 = 3
    29 = 8
    30 = 11
    31 = 9
    32 = 6
    33 = 6
    34 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
